$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Prepare row 10 (brand new row) - copy the bold/centered/bordered
#        style used by column A (A2:A9) onto A10.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. Add the new "Avarage clustering coefficient" property. Row 9's
#        label becomes the new metric, and the property that used to
#        live in row 9 ("The size of largest component") is written out
#        again in the freshly-added row 10.
$ws.Range("B9").Value = "Avarage clustering coefficient"
$ws.Range("B10").Value = "The size of largest component"

# --- 3. Update the numeric values in column C / A per the diff.
$ws.Range("C3").Value = 696
$ws.Range("C4").Value = 1587
$ws.Range("C5").Value = 2.2802
$ws.Range("C6").Value = 0.0066
$ws.Range("C7").Value = 0.0009
$ws.Range("C9").Value = 0.6987

$ws.Range("A10").Value = 8
$ws.Range("C10").Value = 696
